$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("G30").Value = 1.25
$ws.Range("L30").Value = 8
$ws.Range("O30").Value = 1.13
$ws.Range("P30").Value = 6
$ws.Range("Q30").Value = 1.44
$ws.Range("R30").Value = 2.7
$ws.Range("S30").Value = 1.22
$ws.Range("T30").Value = 4
$ws.Range("W30").Value = 9
$ws.Range("Y30").Value = 9
$ws.Range("AA30").Value = 10
$ws.Range("AC30").Value = 21
$ws.Range("AD30").Value = 13
$ws.Range("AF30").Value = 51
$ws.Range("AI30").Value = 41
$ws.Range("AJ30").Value = 23
$ws.Range("AK30").Value = 101
$ws.Range("AL30").Value = 51
$ws.Range("AQ30").Value = 13
$ws.Range("AT30").Value = 4
$ws.Range("AU30").Value = 9
$ws.Range("AZ30").Value = 151
$ws.Range("G31").Value = 1.95
$ws.Range("H31").Value = 3.3
$ws.Range("I31").Value = 4.1
$ws.Range("J31").Value = 2.63
$ws.Range("K31").Value = 2.1
$ws.Range("W31").Value = 6.5
$ws.Range("X31").Value = 9
$ws.Range("Z31").Value = 17
$ws.Range("AA31").Value = 17
$ws.Range("AH31").Value = 11
$ws.Range("AI31").Value = 19
$ws.Range("AJ31").Value = 15
$ws.Range("AO31").Value = 11
$ws.Range("AR31").Value = 51
$ws.Range("AZ31").Value = 81
$ws.Range("G39").Value = 2.88
$ws.Range("H39").Value = 2.8
$ws.Range("I39").Value = 2.55
$ws.Range("Y39").Value = 12
$ws.Range("AC39").Value = 7
$ws.Range("AD39").Value = 5.5
$ws.Range("AI39").Value = 12
$ws.Range("AJ39").Value = 11
$ws.Range("AK39").Value = 26
$ws.Range("AW39").Value = 4.5
